# Fix the "supermaket" typo -> "supermarket" in the bills sheet.
# Cell C3/C4 are shared-formula cells that mirror C2, so fixing the
# source cell (C2) propagates the corrected text to them automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bills")

$ws.Range("C2").Value = "supermarket"

# Move the active selection to C3 (matches the post-edit cursor position).
$ws.Range("C3").Select()
